$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- "MCF" sheet: set several capacity-factor inputs to 1 ---
$wsMcf = $wb.Worksheets.Item("MCF")
$wsMcf.Range("B2").Value = 1
$wsMcf.Range("B3").Value = 1
$wsMcf.Range("B4").Value = 1
$wsMcf.Range("B6").Value = 1
$wsMcf.Range("B10").Value = 1
$wsMcf.Range("B11").Value = 1
$wsMcf.Range("B12").Value = 1
$wsMcf.Range("B13").Value = 1
$wsMcf.Range("B14").Value = 1
$wsMcf.Range("B16").Value = 1

# Move the active selection on the MCF sheet to B17 (matches the saved view state)
$wsMcf.Activate()
$wsMcf.Range("B17").Select()
